$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-05 Thursday" "2026-02-06 Friday"

Replace-Text "88×27=" "57×80="
Replace-Text "79×61=" "15×95="
Replace-Text "42×18=" "92×66="
Replace-Text "75×81=" "84×40="
Replace-Text "27×15=" "69×59="
Replace-Text "41×38=" "60×93="
Replace-Text "72×24=" "43×14="
Replace-Text "16×82=" "21×13="
Replace-Text "17×19=" "43×17="
Replace-Text "47×68=" "54×51="
Replace-Text "63×50=" "41×33="
Replace-Text "97×64=" "60×33="
Replace-Text "11×67=" "39×41="
Replace-Text "95×34=" "30×76="
Replace-Text "77×21=" "36×43="
Replace-Text "61×85=" "75×62="
Replace-Text "25×66=" "82×89="
Replace-Text "19×99=" "47×73="
Replace-Text "68×62=" "63×30="
Replace-Text "69×65=" "13×53="
Replace-Text "76×12=" "48×16="
Replace-Text "42×54=" "86×78="
Replace-Text "24×82=" "60×13="
Replace-Text "23×67=" "60×96="
Replace-Text "20×16=" "89×66="
